# Chinh sua Ke Hoach
# Shifts the plan dates in rows 6-13 one week earlier (the "10/5/2018" text
# marker in F6 becomes "9/28/2018", and every other plan/actual date in
# columns E/F for rows 6-13 moves back by 7 days), then leaves the
# selection on G13 (scrolled so row 4 is the first visible row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("Viet ham Them san pham moi..."): expected-end date text marker.
$ws.Range("E6").Value = 43371
$ws.Range("F6").Value = "                   9/28/2018"

# Row 7 ("Viet ham Xoa san pham..."): actual start/end dates.
$ws.Range("E7").Value = 43372
$ws.Range("F7").Value = 43372

# Row 8 ("Viet ham Mua san pham..."): actual start/end dates.
$ws.Range("E8").Value = 43378
$ws.Range("F8").Value = 43380

# Row 9 ("Kiem tra, sua loi, tinh chinh..."): actual start/end dates.
$ws.Range("E9").Value = 43385
$ws.Range("F9").Value = 43386

# Row 10 ("Len y tuong va tim hieu ve tao giao dien..."): actual start/end dates.
$ws.Range("E10").Value = 43387
$ws.Range("F10").Value = 43387

# Row 11 ("Thiet ke giao dien don gian..."): actual start/end dates.
$ws.Range("E11").Value = 43392
$ws.Range("F11").Value = 43394

# Row 12 ("Kiem tra, sua loi, bo sung va hoan thien..."): actual start/end dates.
$ws.Range("E12").Value = 43399
$ws.Range("F12").Value = 43400

# Row 13 ("Viet bao cao do an..."): actual start/end dates.
$ws.Range("E13").Value = 43401
$ws.Range("F13").Value = 43408

# Leave the selection where the author left it.
[void]$ws.Range("G13").Select()
